$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# ALC row 69
$ws.Range("H69").Value = 22598
$ws.Range("I69").Value = 3400
$ws.Range("J69").Value = 27397.5
$ws.Range("K69").Value = 10200
$ws.Range("L69").Value = 82192.5
$ws.Range("M69").Value = -9326
$ws.Range("N69").Value = -83940.5

# ALC row 72
$ws.Range("H72").Value = 22598
$ws.Range("I72").Value = 3400
$ws.Range("J72").Value = 27397.5
$ws.Range("K72").Value = 30600
$ws.Range("L72").Value = 246577.5
$ws.Range("M72").Value = -26232
$ws.Range("N72").Value = -255313.5

# ALC row 100
$ws.Range("H100").Value = 2884.077
$ws.Range("I100").Value = 1997.6666
$ws.Range("J100").Value = 3150
$ws.Range("K100").Value = 1997.6666
$ws.Range("L100").Value = 3150
$ws.Range("M100").Value = -1456.6666
$ws.Range("N100").Value = -4232

# ALC row 125
$ws.Range("H125").Value = 1240.4286
$ws.Range("I125").Value = 1230.4
$ws.Range("J125").Value = 1265.5
$ws.Range("K125").Value = 11073.6
$ws.Range("L125").Value = 11389.5
$ws.Range("M125").Value = -8613.6
$ws.Range("N125").Value = -16309.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 2483.3333
$ws.Range("I2").Value = 2166.6667
$ws.Range("J2").Value = 2800
$ws.Range("K2").Value = 2166.6667
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = -2053.6667
$ws.Range("N2").Value = -3026

# ARM row 116
$ws.Range("H116").Value = 2483.3333
$ws.Range("I116").Value = 2166.6667
$ws.Range("J116").Value = 2800
$ws.Range("K116").Value = 2166.6667
$ws.Range("L116").Value = 2800
$ws.Range("M116").Value = 127.3332999999998
$ws.Range("N116").Value = -7388

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 2483.3333
$ws.Range("I3").Value = 2166.6667
$ws.Range("J3").Value = 2800
$ws.Range("K3").Value = 2166.6667
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = -2052.6667
$ws.Range("N3").Value = -3028

# BSM row 86
$ws.Range("H86").Value = 3029.5
$ws.Range("I86").Value = 3575
$ws.Range("J86").Value = 2665.8333
$ws.Range("K86").Value = 3575
$ws.Range("L86").Value = 2665.8333
$ws.Range("M86").Value = -2452
$ws.Range("N86").Value = -4911.8333

# BSM row 89
$ws.Range("H89").Value = 3029.5
$ws.Range("I89").Value = 3575
$ws.Range("J89").Value = 2665.8333
$ws.Range("K89").Value = 17875
$ws.Range("L89").Value = 13329.1665
$ws.Range("M89").Value = -12259
$ws.Range("N89").Value = -24561.1665

# BSM row 99
$ws.Range("H99").Value = 1962.8148
$ws.Range("I99").Value = 1591.7858
$ws.Range("J99").Value = 2362.3845
$ws.Range("K99").Value = 1591.7858
$ws.Range("L99").Value = 2362.3845
$ws.Range("M99").Value = -93.78580000000011
$ws.Range("N99").Value = -5358.3845

$ws = $wb.Worksheets.Item("CRP")
# CRP row 99
$ws.Range("H99").Value = 1882.125
$ws.Range("I99").Value = 1646
$ws.Range("K99").Value = 1646
$ws.Range("M99").Value = -148

# CRP row 126
$ws.Range("H126").Value = 1882.125
$ws.Range("I126").Value = 1646
$ws.Range("K126").Value = 4938
$ws.Range("M126").Value = -2468

# CRP row 134
$ws.Range("H134").Value = 1861.3334
$ws.Range("I134").Value = 1813.84
$ws.Range("J134").Value = 2098.8
$ws.Range("K134").Value = 5441.52
$ws.Range("L134").Value = 6296.400000000001
$ws.Range("M134").Value = -2906.52
$ws.Range("N134").Value = -11366.4

$ws = $wb.Worksheets.Item("CUL")
# CUL row 12
$ws.Range("H12").Value = 204.61539
$ws.Range("I12").Value = 253.375
$ws.Range("J12").Value = 182.94444
$ws.Range("K12").Value = 760.125
$ws.Range("L12").Value = 548.83332
$ws.Range("M12").Value = -587.125
$ws.Range("N12").Value = -894.83332

# CUL row 23
$ws.Range("H23").Value = 227.10527
$ws.Range("I23").Value = 85
$ws.Range("J23").Value = 253.75
$ws.Range("K23").Value = 255
$ws.Range("L23").Value = 761.25
$ws.Range("M23").Value = -20
$ws.Range("N23").Value = -1231.25

# CUL row 97
$ws.Range("H97").Value = 508.61703
$ws.Range("I97").Value = 260.05
$ws.Range("J97").Value = 692.7406999999999
$ws.Range("K97").Value = 780.1500000000001
$ws.Range("L97").Value = 2078.2221
$ws.Range("M97").Value = -284.1500000000001
$ws.Range("N97").Value = -3070.2221

# CUL row 131
$ws.Range("H131").Value = 896.9400000000001
$ws.Range("I131").Value = 855.5
$ws.Range("J131").Value = 897.7857
$ws.Range("K131").Value = 2566.5
$ws.Range("L131").Value = 2693.3571
$ws.Range("M131").Value = 2473.5
$ws.Range("N131").Value = -12773.3571

$ws = $wb.Worksheets.Item("LTW")
# LTW row 61
$ws.Range("H61").Value = 2725
$ws.Range("I61").Value = 2160
$ws.Range("J61").Value = 3666.6667
$ws.Range("K61").Value = 2160
$ws.Range("L61").Value = 3666.6667
$ws.Range("M61").Value = -1958
$ws.Range("N61").Value = -4070.6667

# LTW row 82
$ws.Range("H82").Value = 2637.3572
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2637.3572
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2637.3572
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3359.3572

# LTW row 85
$ws.Range("H85").Value = 2637.3572
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2637.3572
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2637.3572
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5133.3572

# LTW row 113
$ws.Range("H113").Value = 2725
$ws.Range("I113").Value = 2160
$ws.Range("J113").Value = 3666.6667
$ws.Range("K113").Value = 2160
$ws.Range("L113").Value = 3666.6667
$ws.Range("M113").Value = 10
$ws.Range("N113").Value = -8006.6667

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 3333.3333
$ws.Range("I62").Value = 3333.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3333.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2709.3333
$ws.Range("N62").ClearContents()

# WVR row 65
$ws.Range("H65").Value = 3333.3333
$ws.Range("I65").Value = 3333.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16666.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13546.6665
$ws.Range("N65").ClearContents()

# WVR row 81
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2939
$ws.Range("N81").ClearContents()

# WVR row 84
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -14696
$ws.Range("N84").ClearContents()
